$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.305.37"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.605.82"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "543.59"
$ws.Range("E5").Value = "  +4.18%  "
$ws.Range("D6").Value = "140.93"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "6.46"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").Value = "0.333"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "3.065.80"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "59.228.59"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "20.53"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.586.02"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "342.93"
$ws.Range("D19").Value = "4.36"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "67.52"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "0.408"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "7.21"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("E30").Value = "  +8.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.80"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "149.39"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.10"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "0.833"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "277.18"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").Value = "10.74"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "0.0955"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").Value = "1.942.83"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "18.36"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.90"
$ws.Range("E51").Value = "  -2.80%  "
